# feat: correções nas tabelas fato e dimensão a serem utilizados no desafio técnico
#
# The "Lideres" (inspector name) column on the dim_insp sheet had trailing
# whitespace on every entry. This trims that trailing whitespace from each
# name in column B (data rows 2-43), leaving ID / Base / Setor untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dim_insp")

# Find the last used row in column A (ID) so we cover every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 43 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value2
    if ($current -ne $null) {
        $trimmed = $current.TrimEnd()
        if ($trimmed -ne $current) {
            $cell.Value2 = $trimmed
        }
    }
}

# Leave the selection where the user last clicked while cleaning up the list.
$ws.Range("B6").Select()
